# Apply the "GET functions integration" update to the amazon product data
# workbook:
#   - drop the duplicate "T-Shirts" sheet (ID_dda333c) entirely
#   - rename the "SSD" sheet (ID_ca64249) to ID_03f327c and refresh its
#     data to the "Mens Cotton Jacket" product
#   - de-duplicate the remaining backpack sheet (ID_7e0f17b), keeping a
#     single data row, and bump its tracked date

$excel.DisplayAlerts = $false | Out-Null
$wb = $excel.ActiveWorkbook

# --- Sheet "ID_7e0f17b": remove the duplicated row 3, refresh the date ---
$wsBackpack = $wb.Worksheets.Item("ID_7e0f17b")
$wsBackpack.Range("B2").Value = "17/02/2025"
[void]$wsBackpack.Rows.Item(3).Delete()

# --- Remove the duplicate T-Shirts sheet entirely ---
$wsShirt = $wb.Worksheets.Item("ID_dda333c")
[void]$wsShirt.Delete()

# --- Rename the SSD sheet and replace its tracked product with the jacket ---
$wsJacket = $wb.Worksheets.Item("ID_ca64249")
$wsJacket.Name = "ID_03f327c"

$wsJacket.Range("A2").Value = 55.99
$wsJacket.Range("B2").Value = "17/02/2025"
$wsJacket.Range("C2").Value = "Mens Cotton Jacket"
$wsJacket.Range("D2").Value = "great outerwear jackets for Spring/Autumn/Winter, suitable for many occasions, such as working, hiking, camping, mountain/rock climbing, cycling, traveling or other outdoors. Good gift choice for you or your family member. A warm hearted love to Father, husband or son in this thanksgiving or Christmas Day."
$wsJacket.Range("E2").Value = "ID_03f327c"
